# "added changings in write cast logic"
# Rework the RestAssuredTest sheet: drop the old "Condititions/Assertions"
# title row, promote the former header row to row 1 with renamed/shortened
# column captions, widen the columns to fit the new captions, and refresh
# the hyperlinks so they keep pointing at the (now shifted) URL cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RestAssuredTest")

# Remove the merged title row (old row 1, "Condititions/Assertions").
# This shifts every row up by one: old headers (row2) -> row1,
# old data rows (3,4,5) -> rows 2,3,4, old blank rows (6-10) -> rows 5-9.
$ws.Range("D1:G1").UnMerge()
$ws.Rows.Item(1).Delete()

# Rename / shorten the header captions on the new row 1.
$ws.Range("A1").Value = "API"
$ws.Range("B1").Value = "Method"
$ws.Range("C1").Value = "Request Payload "
$ws.Range("D1").Value = "Response Contains"
$ws.Range("G1").Value = "Response Size"
$ws.Range("F1").Value = "Response time"
$ws.Range("E1").Value = "Response Code"

# New column widths to accommodate the longer header text.
$ws.Columns.Item(1).ColumnWidth = 21.8
$ws.Columns.Item(4).ColumnWidth = 28.0833333333333
$ws.Columns.Item(5).ColumnWidth = 20.6
$ws.Columns.Item(6).ColumnWidth = 38.26
$ws.Columns.Item(7).ColumnWidth = 48.6

# The header row now wraps onto fewer lines given the wider columns.
$ws.Rows.Item(1).RowHeight = 31.3

# Hyperlinks need to be re-anchored: they used to sit on A3/A4, but after
# removing the title row the same URLs now live on A2/A3.
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://emumba-test.herokuapp.com/trend")
$ws.Hyperlinks.Add($ws.Range("A3"), "https://emumba-test.herokuapp.com/user")

# Update the active selection to match the saved view state.
$ws.Range("E3").Select() | Out-Null
